# Tue, Aug 04, 2020 12:05:05 AM
#
# The deck's slide master is wired to the "Integral" theme
# (dk1=000000 lt1=FFFFFF dk2=455F51 lt2=E3DED1 accent1=99CB38 accent2=63A537
#  accent3=E6D024 accent4=CC9700 accent5=4EB3CF accent6=378DA6 hlink=6B9F25
#  folHlink=B26B02) while a second, unused theme part still carries the
# original default "Office Theme" palette. This edit swaps the two so the
# presentation's live design reverts to the stock Office colour scheme
# (dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6 accent1=5B9BD5 accent2=ED7D31
#  accent3=A5A5A5 accent4=FFC000 accent5=4472C4 accent6=70AD47 hlink=0563C1
#  folHlink=954F72).
#
# PowerPoint's theme colours are edited through
# ThemeColorScheme.Colors(i).RGB (i = 1..12, in the fixed order
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) -- RGB() packs as
# R + G*256 + B*65536, the usual VBA colour encoding.

function RGBColor([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# The live/active theme (backing the slide master every slide inherits
# from) is reachable via SlideMaster.Theme.
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Colors(1).RGB  = RGBColor 0x00 0x00 0x00   # dk1      000000
$tcs.Colors(2).RGB  = RGBColor 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Colors(3).RGB  = RGBColor 0x44 0x54 0x6A   # dk2      44546A
$tcs.Colors(4).RGB  = RGBColor 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Colors(5).RGB  = RGBColor 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Colors(6).RGB  = RGBColor 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Colors(7).RGB  = RGBColor 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Colors(8).RGB  = RGBColor 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Colors(9).RGB  = RGBColor 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Colors(10).RGB = RGBColor 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Colors(11).RGB = RGBColor 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Colors(12).RGB = RGBColor 0x95 0x4F 0x72   # folHlink 954F72
